# dnn pool cidr change
#
# Slide 1 (sldId 256 / slide1.xml) edits:
#   1. Shrink the "Rectángulo: esquinas redondeadas 2" dashed box (id 3)
#      so it no longer spans the full width (its right edge now leaves
#      room for a new box to its right).
#   2. Slide the "2" icon + its numeral label (ids 35 / 36) from inside
#      that box out to the right, near the new CIDR box.
#   3. Add a brand-new dashed "no fill" rounded-rectangle placeholder
#      box (a duplicate of shape id 3's style) to the right of the
#      existing two, to hold the new DNN pool CIDR block.

function Get-ShapeById {
    param($Shapes, $Id)
    foreach ($sh in $Shapes) {
        if ($sh.Id -eq $Id) { return $sh }
    }
    return $null
}

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)

# ---------------------------------------------------------------
# 1) Shrink shape id=3 ("Rectángulo: esquinas redondeadas 2").
#    ext cx: 4658050 EMU -> 3060045 EMU (height/position unchanged).
# ---------------------------------------------------------------
$box2 = Get-ShapeById $slide.Shapes 3
$box2.Width = 240.9484251968504

# ---------------------------------------------------------------
# 2) Reposition the "2" icon picture (id=35) and its numeral
#    textbox (id=36) out to the right of the shrunk box.
# ---------------------------------------------------------------
$icon2 = Get-ShapeById $slide.Shapes 35
$icon2.Left = 600.7371653543307
$icon2.Top  = 197.9015748031496

$label2 = Get-ShapeById $slide.Shapes 36
$label2.Left = 602.6133170866142
$label2.Top  = 197.9015748031496

# ---------------------------------------------------------------
# 3) Add the new dashed rounded-rectangle box for the DNN pool CIDR.
#    The cleanest way to reproduce the exact same style/txBody as the
#    existing boxes (no fill, dashed 2pt accent1 outline, style refs,
#    centered empty paragraph) is to duplicate shape id=3 and then
#    move/resize/rename the copy, rather than rebuild it via
#    Shapes.AddShape (which creates a plain, unstyled shape).
#
#    This deck's shape-id allocator hands out the *lowest unused* id
#    rather than always incrementing a running counter, so id 47 is
#    the 6th gap available (17, 31, 33, 37, 45, 47). We duplicate six
#    times and discard the first five throwaway copies, keeping only
#    the 6th (id 47) as the real new shape.
# ---------------------------------------------------------------
$spawned = @()
for ($i = 0; $i -lt 6; $i++) {
    $dup = $box2.Duplicate()
    $spawned += $dup.Item(1)
}

for ($i = 0; $i -lt ($spawned.Count - 1); $i++) {
    $spawned[$i].Delete()
}

$newBox = $spawned[$spawned.Count - 1]
$newBox.Name = "Rectángulo: esquinas redondeadas 46"
$newBox.Left   = 439.23614173228344
$newBox.Top    = 197.05426196850394
$newBox.Width  = 194.56701787401576
$newBox.Height = 103.95371078740158
